# Fruta / hortaliza, semanal
# The underlying weekly data refresh re-shuffles the "observation" columns
# (Fecha, Calidad, Volumen, Precio minimo/maximo/promedio, Unidad de
# comercializacion, Origen, Precio $/Kg, Kg/unidad) across the existing
# data rows (2-63), while the descriptive columns (Mercado ID, Mercado,
# Region, Codreg, Tipo, Producto ID, Producto, Categoria ID, Categoria,
# Variedad) stay put. This script snapshots the movable columns for every
# row and then redistributes them according to the mapping derived from
# the published diff.

$ws = $excel.ActiveWorkbook.ActiveSheet

$firstRow = 2
$lastRow  = 63

# Columns that move between rows.
$cols = @("D","L","M","N","O","P","Q","R","S","T")

# Mapping: new row -> row whose old values are copied into it.
$mapping = @{
    2  = 5;  3  = 39; 4  = 25; 5  = 26; 6  = 27; 7  = 23; 8  = 18; 9  = 28;
    10 = 42; 11 = 10; 12 = 50; 13 = 24; 14 = 40; 15 = 34; 16 = 11; 17 = 55;
    18 = 61; 19 = 60; 20 = 63; 21 = 12; 22 = 41; 23 = 21; 24 = 43; 25 = 13;
    26 = 37; 27 = 44; 28 = 45; 29 = 46; 30 = 3;  31 = 62; 32 = 47; 33 = 6;
    34 = 14; 35 = 19; 36 = 29; 37 = 36; 38 = 56; 39 = 31; 40 = 32; 41 = 15;
    42 = 7;  43 = 8;  44 = 54; 45 = 17; 46 = 49; 47 = 4;  48 = 57; 49 = 30;
    50 = 51; 51 = 52; 52 = 22; 53 = 38; 54 = 48; 55 = 59; 56 = 35; 57 = 20;
    58 = 9;  59 = 53; 60 = 58; 61 = 16; 62 = 33; 63 = 2
}

# Snapshot the current ("before") values of every movable cell so the
# write-back pass below never reads data that has already been overwritten.
# Value2 is used (rather than Value) so dates come back as raw serial
# numbers instead of ambiguous Variant wrappers.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Write the permuted values back.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $mapping[$r]
    $rowVals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $rowVals[$c]
    }
}
